# The "messages" sheet gained 9 new chat-log rows (rows 5-13) and row 4's
# IsEdited/IsDeleted flags were corrected from booleans to numeric 0/1 to
# match the export format used by the newly appended rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 4 (existing row): IsEdited/IsDeleted -> numeric ---
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1

# --- Row 5: message #4 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "U001"
$ws.Range("C5").Value = "آقای گلستانی"
$ws.Range("D5").Value = "سلام خدمت همکاران محترم فروش"
$ws.Range("E5").Value = "2026-02-03 23:03:21"
$ws.Range("F5").Value = "'1404/11/14"
$ws.Range("G5").Value = "23:03"
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = "none"
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""

# --- Row 6: message #5 ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "U001"
$ws.Range("C6").Value = "آقای گلستانی"
$ws.Range("D6").Value = "سلام"
$ws.Range("E6").Value = "2026-02-03 23:03:33"
$ws.Range("F6").Value = "'1404/11/14"
$ws.Range("G6").Value = "23:03"
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = "none"
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = ""

# --- Row 7: message #6 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "U001"
$ws.Range("C7").Value = "آقای گلستانی"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "2026-02-03 23:03:54"
$ws.Range("F7").Value = "'1404/11/14"
$ws.Range("G7").Value = "23:03"
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = "image"
$ws.Range("M7").Value = "1029.jpg"
$ws.Range("N7").Value = "static/chat_uploads\20260203_230354_0931990f_1029.jpg"
$ws.Range("O7").Value = "109.5 KB"

# --- Row 8: message #7 ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "U001"
$ws.Range("C8").Value = "آقای گلستانی"
$ws.Range("D8").Value = "شسیشسی"
$ws.Range("E8").Value = "2026-02-03 23:04:08"
$ws.Range("F8").Value = "'1404/11/14"
$ws.Range("G8").Value = "23:04"
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = "none"
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("O8").Value = ""

# --- Row 9: message #8 ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "U001"
$ws.Range("C9").Value = "آقای گلستانی"
$ws.Range("D9").Value = "ثقصثقصثقصث"
$ws.Range("E9").Value = "2026-02-03 23:05:03"
$ws.Range("F9").Value = "'1404/11/14"
$ws.Range("G9").Value = "23:05"
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = "none"
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("O9").Value = ""

# --- Row 10: message #9 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "U001"
$ws.Range("C10").Value = "آقای گلستانی"
$ws.Range("D10").Value = "شسیشسی"
$ws.Range("E10").Value = "2026-02-03 23:05:22"
$ws.Range("F10").Value = "'1404/11/14"
$ws.Range("G10").Value = "23:05"
$ws.Range("H10").Value = $false
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = "none"
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = ""
$ws.Range("O10").Value = ""

# --- Row 11: message #10 ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "U001"
$ws.Range("C11").Value = "آقای گلستانی"
$ws.Range("D11").Value = "سلاک"
$ws.Range("E11").Value = "2026-02-03 23:11:36"
$ws.Range("F11").Value = "'1404/11/14"
$ws.Range("G11").Value = "23:11"
$ws.Range("H11").Value = $false
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = "none"
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("O11").Value = ""

# --- Row 12: message #11 ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "U001"
$ws.Range("C12").Value = "آقای گلستانی"
$ws.Range("D12").Value = "سلام"
$ws.Range("E12").Value = "2026-02-03 23:11:44"
$ws.Range("F12").Value = "'1404/11/14"
$ws.Range("G12").Value = "23:11"
$ws.Range("H12").Value = $false
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = "none"
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = ""

# --- Row 13: message #12 ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "U001"
$ws.Range("C13").Value = "آقای گلستانی"
$ws.Range("D13").Value = "سلام خدمت همکاران محترم فروش این نرم افزاربرای کمک به شما جهت شناخت مسیر خودتون طراحی شده"
$ws.Range("E13").Value = "2026-02-04 08:29:27"
$ws.Range("F13").Value = "'1404/11/15"
$ws.Range("G13").Value = "08:29"
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = $false
$ws.Range("K13").Value = $false
$ws.Range("L13").Value = "none"
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("O13").Value = ""
